# Add six new modeling-parameter columns (log_lik, al_win, al_loss, kappa,
# beta_win, beta_loss) to the existing Table3 on "Sheet 1", fill in the
# header names and the 28 rows of per-subject values, and restore the
# selection to match the post-edit workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- grow the table by six columns (BQ:BV) ---------------------------------
for ($i = 0; $i -lt 6; $i++) {
    $lo.ListColumns.Add() | Out-Null
}

# Header row. Setting the header cell's Value is what drives the
# corresponding ListColumn's Name, so just write the cells directly.
# (Column order on the sheet is log_lik, al_win, al_loss, kappa, beta_win,
#  beta_loss; header cells are written in this left-to-right order except
# that "kappa" is entered before "al_loss" to reproduce the original
# shared-string insertion order.)
$ws.Range("BQ1").Value = "log_lik"
$ws.Range("BR1").Value = "al_win"
$ws.Range("BT1").Value = "kappa"
$ws.Range("BS1").Value = "al_loss"
$ws.Range("BU1").Value = "beta_win"
$ws.Range("BV1").Value = "beta_loss"

# --- per-subject modeling output (rows 2-29) -------------------------------
# Columns, in order: log_lik, al_win, al_loss, kappa, beta_win, beta_loss
$data = @(
    @(-195.24523001490601,0.27650000000000002,0.74339999999999995,-0.59751766350473201,0.70918673554329703,-0.42291431385054601),
    @(-165.26080230379301,0.32290000000000002,0.88219999999999998,-0.33299572561956597,1.02957655196939,0.46879038000710099),
    @(-101.47496579892901,0.23,0.92769999999999997,-1.76299798096754,1.70106939969715,1.3525657350028),
    @(-65.755395678022396,0.21859999999999999,0.61099999999999999,-0.99090813395190702,2.01221574468331,1.4327477919953999),
    @(-193.48568947408901,0.24859999999999999,0.92749999999999999,-0.29302557277204,0.645611501496689,0.45243428769427402),
    @(-113.570517310347,0.2039,0.49730000000000002,-0.95623048809300604,1.6559707678286,1.6041850537107101),
    @(-81.0274810010431,0.29530000000000001,0.76160000000000005,-1.22511215610414,1.8751071144323701,1.3270453347473901),
    @(-153.28737607999301,0.37509999999999999,0.87639999999999996,0.155182060518574,1.0768682501890099,0.050772242594244003),
    @(-36.2772334443777,0.1147,0.6784,-3.4697156398402198,2.7614574858654,2.3663676394868398),
    @(-176.20970671031299,0.30130000000000001,0.79469999999999996,-2.2695774014510102,0.91615615799420203,0.72365446505258002),
    @(-129.59821727707401,0.47770000000000001,0.90669999999999995,-0.31101176545220599,1.1192423119757899,1.03779990004914),
    @(-60.571385672718101,0.31480000000000002,0.61270000000000002,-1.32791011237045,2.0580095457991399,1.43560488757254),
    @(-94.976580414405802,0.42920000000000003,0.92079999999999995,-0.43466934283172798,1.5280735742400799,0.80703860568608998),
    @(-105.125024479212,0.41760000000000003,0.82040000000000002,-0.73835003836202995,1.6033147288251199,0.082613650931698895),
    @(-57.811669798091401,0.30570000000000003,0.76500000000000001,-1.47391208041934,2.09398730181656,1.75676989537803),
    @(-115.194053235957,0.26200000000000001,0.82210000000000005,-1.0707873426015899,1.5906079169978899,0.69118923646263297),
    @(-100.150514235363,0.31340000000000001,0.89759999999999995,0.369544565129868,1.5011309256635199,0.99966629456970302),
    @(-176.71758865828599,0.1067,0.72709999999999997,-1.2759578059546901,1.4511630742980901,-0.054760836326000199),
    @(-101.67088325247801,0.317,0.81130000000000002,-0.54412988590080202,1.7790091285483101,1.12841757442276),
    @(-134.796221624216,0.36940000000000001,0.89710000000000001,-1.1831840361529,1.3234407918386299,0.482531567602069),
    @(-98.593096553727406,0.26889999999999997,0.87370000000000003,-1.2795619426334099,1.8796881403993599,0.82152876434571698),
    @(-117.579142180875,0.32040000000000002,0.86150000000000004,-0.71704458688894102,1.4631170900794599,1.08708767510257),
    @(-162.18186632209,0.12759999999999999,0.25459999999999999,-2.0520179673243502,1.5873608032834201,0.96186727280885898),
    @(-58.184112111257299,0.30249999999999999,0.50680000000000003,-1.0714361175579501,1.8317257290739599,1.8052383920493),
    @(-75.349284475398306,0.33450000000000002,0.32969999999999999,-1.5480920879250499,1.71137806116449,2.2558751583223802),
    @(-114.67111241329999,0.32540000000000002,0.89729999999999999,-0.81103065015539,1.70783692537763,-0.125085432174869),
    @(-49.3868939940239,0.27060000000000001,0.72709999999999997,-2.0177914044123502,2.4039996018230299,1.8365387948684899),
    @(-82.294699015430794,0.19639999999999999,0.23749999999999999,-0.36222088391919599,1.96481284620143,1.6009380656563701)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($rowNum, 69 + $j).Value = $vals[$j]
    }
}

# --- restore the on-screen selection to the cell active after the edit ----
$ws.Range("BO11").Select()
